$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-10 Thursday", "2025-04-11 Friday"),
    @("32×36=", "27×43="),
    @("11×58=", "49×57="),
    @("32×67=", "98×72="),
    @("95×46=", "53×45="),
    @("70×90=", "80×55="),
    @("74×34=", "60×84="),
    @("52×62=", "89×48="),
    @("11×56=", "33×99="),
    @("31×53=", "65×27="),
    @("86×44=", "96×21="),
    @("21×18=", "53×77="),
    @("73×62=", "13×72="),
    @("59×20=", "14×75="),
    @("89×63=", "56×56="),
    @("21×60=", "46×26="),
    @("43×18=", "28×38="),
    @("74×68=", "52×37="),
    @("44×58=", "30×79="),
    @("95×67=", "98×96="),
    @("86×55=", "53×39="),
    @("32×73=", "99×87="),
    @("83×79=", "53×72="),
    @("39×65=", "39×88="),
    @("37×91=", "98×39="),
    @("22×44=", "91×87=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
